$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("0.0703539177775383","0.9795734882354736","0.005070631392300129","0.9990717172622681"),
    @("0.0113139022141695","0.9982664585113525","0.003721604123711586","0.9993369579315186"),
    @("0.00702159246429801","0.9986215233802795","0.001183720538392663","0.9993369579315186"),
    @("0.003760518273338675","0.9991227984428406","0.0002141550794476643","1"),
    @("0.002411472843959928","0.9994778633117676","9.648475679568946E-05","1"),
    @("0.001630087848752737","0.9995405077934265","3.223944440833293E-05","1"),
    @("0.002014506608247757","0.9995822906494141","2.508226680220105E-05","1"),
    @("0.00104352948255837","0.9997493624687195","4.614570934791118E-05","1"),
    @("0.001174238394014537","0.9996867179870605","0.0001710880460450426","0.9999337196350098"),
    @("0.000663137121591717","0.9998747110366821","4.064757740707137E-05","1"),
    @("0.001226450898684561","0.9998329281806946","1.332003193965647E-05","1"),
    @("0.0003111085679847747","0.9999164342880249","4.679050562117482E-06","1"),
    @("0.0007372401305474341","0.9997702836990356","1.64404127644957E-06","1"),
    @("0.0003883748140651733","0.9998955726623535","1.090363184630405E-06","1"),
    @("0.0007096432964317501","0.9998329281806946","4.650347364076879E-06","1"),
    @("0.0003715037018992007","0.9998955726623535","6.232102691683394E-07","1"),
    @("0.000333794770995155","0.9999164342880249","1.69585788967197E-07","1"),
    @("0.0007025942904874682","0.999791145324707","7.282476872205734E-06","1"),
    @("0.0002376369229750708","0.9999164342880249","5.816940529257408E-07","1"),
    @("0.0004891667049378157","0.9998955726623535","2.341200513455988E-07","1"),
    @("0.0005570605862885714","0.999791145324707","7.617665573889099E-07","1"),
    @("2.429872984066606E-05","1","3.120296980796411E-07","1"),
    @("0.000293505028821528","0.9998955726623535","4.836655307371984E-07","1"),
    @("0.0005031879409216344","0.9998747110366821","0.0001252242218470201","0.9999337196350098"),
    @("6.838005356257781E-05","0.9999582171440125","2.482701404460386E-07","1"),
    @("3.947609366150573E-05","0.9999791383743286","4.202100711836465E-08","1"),
    @("0.0001154251804109663","0.9999791383743286","7.362027787394254E-08","1"),
    @("0.0003693216713145375","0.9998955726623535","4.527093011574834E-08","1"),
    @("0.0002016418147832155","0.9999164342880249","7.234248045051572E-08","1"),
    @("0.0001641871349420398","0.9999791383743286","8.22650861209695E-07","1"),
    @("5.290194894769229E-05","0.9999791383743286","1.711826946859674E-08","1"),
    @("0.0009313809568993747","0.9999164342880249","1.256234190805117E-05","1"),
    @("0.000539662956725806","0.9998747110366821","3.499355116787228E-08","1"),
    @("0.000206370372325182","0.9999373555183411","4.434132616637498E-09","1"),
    @("0.0004632231430150568","0.9999791383743286","1.039889028220387E-07","1"),
    @("0.0002641478786244988","0.9999373555183411","3.613794419266014E-08","1"),
    @("0.0004862137429881841","0.9999373555183411","2.187617731408409E-08","1"),
    @("0.0001101760863093659","0.9999791383743286","6.789515172300753E-09","1"),
    @("7.509037095587701E-05","0.9999582171440125","3.881399379679351E-08","1"),
    @("0.0005306308157742023","0.9999373555183411","1.312072894776861E-09","1"),
    @("9.13270196178928E-05","0.9999582171440125","1.187619318443467E-06","1"),
    @("0.0003492504765745252","0.9999582171440125","1.738888255609083E-09","1"),
    @("0.0006606105016544461","0.9998747110366821","3.232746514925111E-09","1"),
    @("0.0002859742671716958","0.9999373555183411","2.529296994069341E-09","1"),
    @("0.0002667694934643805","0.9999791383743286","1.201417854090892E-09","1"),
    @("0.0005186922498978674","0.9998955726623535","8.898423971004377E-07","1"),
    @("4.826203803531826E-05","0.9999791383743286","6.283572773213564E-09","1"),
    @("1.567370054544881E-05","1","4.523637286979465E-08","1"),
    @("7.722579539404251E-06","1","7.516420374997779E-09","1"),
    @("6.004908300383249E-06","1","7.113668681890672E-11","1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = [double]$row[0]
    $ws.Cells.Item($r, 2).Value = [double]$row[1]
    $ws.Cells.Item($r, 3).Value = [double]$row[2]
    $ws.Cells.Item($r, 4).Value = [double]$row[3]
}
